# #5: property boat&car done
# Fix the "汽車" (vehicle) sheet: the header row had leaked data values instead
# of proper column labels, and several metadata columns (property_category,
# category, date, legislator_name, legislator_id, source_file, index) that
# already exist on the other sheets were missing here. Also add the new
# "capacity" field (engine displacement) to the header.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)   # 汽車 (car) sheet

# ---- Row 1: header labels (B1:G1 already exist with the bold/bordered style) ----
$ws.Cells.Item(1,2).Value = "name"
$ws.Cells.Item(1,3).Value = "capacity"
$ws.Cells.Item(1,4).Value = "owner"
$ws.Cells.Item(1,5).Value = "register_date"
$ws.Cells.Item(1,6).Value = "register_reason"
$ws.Cells.Item(1,7).Value = "acquire_value"

# New header cells H1:N1 - copy the existing header style from G1 first so the
# new cells pick up the same bold/centered/bordered formatting, then set text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Cells.Item(1,8).Value = "property_category"

$ws.Range("G1").Copy($ws.Range("I1"))
$ws.Cells.Item(1,9).Value = "category"

$ws.Range("G1").Copy($ws.Range("J1"))
$ws.Cells.Item(1,10).Value = "date"

$ws.Range("G1").Copy($ws.Range("K1"))
$ws.Cells.Item(1,11).Value = "legislator_name"

$ws.Range("G1").Copy($ws.Range("L1"))
$ws.Cells.Item(1,12).Value = "legislator_id"

$ws.Range("G1").Copy($ws.Range("M1"))
$ws.Cells.Item(1,13).Value = "source_file"

$ws.Range("G1").Copy($ws.Range("N1"))
$ws.Cells.Item(1,14).Value = "index"

# ---- Row 2: data for 曰產Livina (keep existing B2:G2 values, add H2:N2) ----
$ws.Range("G2").Copy($ws.Range("H2"))
$ws.Cells.Item(2,8).Value = "land"

$ws.Range("G2").Copy($ws.Range("I2"))
$ws.Cells.Item(2,9).Value = "normal"

$ws.Range("G2").Copy($ws.Range("J2"))
$ws.Cells.Item(2,10).Value = "2011-11-22"

$ws.Range("G2").Copy($ws.Range("K2"))
$ws.Cells.Item(2,11).Value = "邱文彥"

$ws.Range("G2").Copy($ws.Range("L2"))
$ws.Cells.Item(2,12).Value = 1743

$ws.Range("G2").Copy($ws.Range("M2"))
$ws.Cells.Item(2,13).Value = "tmpf3df1"

$ws.Range("G2").Copy($ws.Range("N2"))
$ws.Cells.Item(2,14).Value = 35

# ---- Row 3: data for 日產Sentra (keep existing B3:G3 values, add H3:N3) ----
$ws.Range("G3").Copy($ws.Range("H3"))
$ws.Cells.Item(3,8).Value = "land"

$ws.Range("G3").Copy($ws.Range("I3"))
$ws.Cells.Item(3,9).Value = "normal"

$ws.Range("G3").Copy($ws.Range("J3"))
$ws.Cells.Item(3,10).Value = "2011-11-22"

$ws.Range("G3").Copy($ws.Range("K3"))
$ws.Cells.Item(3,11).Value = "邱文彥"

$ws.Range("G3").Copy($ws.Range("L3"))
$ws.Cells.Item(3,12).Value = 1743

$ws.Range("G3").Copy($ws.Range("M3"))
$ws.Cells.Item(3,13).Value = "tmpf3df1"

$ws.Range("G3").Copy($ws.Range("N3"))
$ws.Cells.Item(3,14).Value = 36

Write-Output "done"
